# Daily automated update: append the 16/5/2020 row to the
# "Condicion_Pacientes" table on Hoja1 (Pruebas Realizadas / Pruebas
# Positivas reported for that date; the remaining columns are not yet
# available and stay blank, same as the previous day's pattern).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item("Condicion_Pacientes")

# Expanding the table automatically grows its range/autoFilter and the
# sheet dimension, mirroring what Excel does when a new record is typed
# into the row right below the table.
$newRow = $lo.ListRows.Add()
$r = $newRow.Range

$r.Cells.Item(1, 1).Value = "16/5/2020"
$r.Cells.Item(1, 2).Value = 1397
$r.Cells.Item(1, 3).Value = 120

# Match the right-aligned date style used by the other text-based date
# entries further up the column (e.g. A46/A47).
$r.Cells.Item(1, 1).HorizontalAlignment = -4152

# Leave the selection on the next column of the new row, same as after
# manually tabbing through entering the row's data.
[void]$r.Cells.Item(1, 4).Select()
